# feat(employee): import employee with relation
#
# Adds 5 new "relation id" columns (P:T) to the employee import template:
#   P = ID Kantor/Cabang
#   Q = ID Divisi
#   R = ID Posisi Jabatan
#   S = ID Level Jabatan
#   T = ID Status Karyawan
#
# New shared-string entries must land in this exact order so the
# sharedStrings table matches the target (append order = first-use order):
#   ID Posisi Jabatan, ID Level Jabatan, ID Status Karyawan,
#   ID Kantor/Cabang, ID Divisi

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write header cells in the order that makes the shared-string table's
# new entries line up: R, S, T first, then P, Q.
$ws.Range("R1").Value = "ID Posisi Jabatan"
$ws.Range("S1").Value = "ID Level Jabatan"
$ws.Range("T1").Value = "ID Status Karyawan"
$ws.Range("P1").Value = "ID Kantor/Cabang"
$ws.Range("Q1").Value = "ID Divisi"

# Column widths for the new columns (best-fit-ish, matching the authored
# template as closely as this engine's width quantization allows).
$ws.Columns(16).ColumnWidth = 15.736979166666666  # P -> ~16.57
$ws.Columns(17).ColumnWidth = 15.736979166666666  # Q -> ~16.57
$ws.Columns(18).ColumnWidth = 12.592447916666666  # R -> ~13.43
$ws.Columns(19).ColumnWidth = 12.022135416666666  # S -> ~12.86
$ws.Columns(20).ColumnWidth = 14.736979166666666  # T -> ~15.57

# Move the view so the newly added columns are visible and select Q10,
# matching the author's saved cursor/scroll position.
$win = $excel.ActiveWindow
$win.ScrollColumn = 10
$win.ScrollRow = 1
$ws.Range("Q10").Select() | Out-Null
